# addMitra.xlsx template update: "perkondisian dan seeder"
#
# - Inserts a new "kode_kabupaten" column ahead of kode_kecamatan, which
#   pushes the old kode_desa header one column to the right (into the slot
#   formerly occupied by jenis_kelamin) and jenis_kelamin one column right
#   again (into the slot formerly occupied by status_pekerjaan, which is
#   retired - detail_pekerjaan now lives there instead).
# - Updates the matching example/seed values in row 2 to match.
# - Moves the green "required" header highlight off kode_desa and onto the
#   relocated jenis_kelamin column; detail_pekerjaan picks up the yellow
#   "optional" highlight that used to only mark the trailing date columns.
# - Deletes the now-stale legend textbox that explained the old
#   jenis_kelamin/status_pekerjaan coding scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
$ws.Range("D1").Value = "kode_kabupaten"
$ws.Range("F1").Value = "kode_desa"
$ws.Range("G1").Value = "jenis_kelamin"

# --- Example data row (row 2) --------------------------------------------
$ws.Range("D2").Value = "16"
$ws.Range("F2").Value = "001"
$ws.Range("G2").Value = "Pr"

# --- Header highlight colours --------------------------------------------
# F1 (now kode_desa) should be plain again - copy formatting from a
# neighbouring plain header cell.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# H1 (now detail_pekerjaan) picks up the yellow "optional" highlight that
# the trailing date columns (K1/L1) already use.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Remove the obsolete jenis_kelamin / status_pekerjaan legend textbox -
$ws.Shapes.Item("TextBox 2").Delete()

# --- Minor view-state touch-ups (zoom / last selection) ------------------
$excel.ActiveWindow.Zoom = 87
$ws.Range("G7").Select() | Out-Null
